$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 144, shifting all existing
# records (old rows 144-245) down by one (new rows 145-246).
$ws.Rows("144:144").Insert()

# Populate the newly inserted row 144 with the new weekly record.
# Columns A,B,C,E,F,G,H,I,J,K,Q,R,T carry over the same values as the
# record that used to sit in row 144 (now row 145); only D, L, M, N, O,
# P and S differ for the new entry.
$ws.Range("A144").Value = 11
$ws.Range("B144").Value = "Vega Monumental Concepción"
$ws.Range("C144").Value = "Bíobío"
$ws.Range("D144").Value = 44978
$ws.Range("E144").Value = 8
$ws.Range("F144").Value = "Fruta"
$ws.Range("G144").Value = 100108
$ws.Range("H144").Value = "Tropicales y subtropicales"
$ws.Range("I144").Value = 100108005
$ws.Range("J144").Value = "Piña"
$ws.Range("K144").Value = "Caramelo"
$ws.Range("L144").Value = "Primera"
$ws.Range("M144").Value = 270
$ws.Range("N144").Value = 23000
$ws.Range("O144").Value = 24000
$ws.Range("P144").Value = 23556
$ws.Range("Q144").Value = "$/caja 14 unidades"
$ws.Range("R144").Value = "Ecuador"
$ws.Range("S144").Value = 1683
$ws.Range("T144").Value = 14
